$wb = $excel.ActiveWorkbook

# Suppress the "delete sheet" confirmation dialog, if supported
$excel.DisplayAlerts = $false

# Rename the first sheet to "Parts List"
$wb.Worksheets.Item("Sheet1").Name = "Parts List"

# Remove the now-unused Sheet2 and Sheet3
$null = $wb.Worksheets.Item("Sheet2").Delete()
$null = $wb.Worksheets.Item("Sheet3").Delete()

$excel.DisplayAlerts = $true
